# Fruta / hortaliza, semanal
#
# A new daily price record is inserted at row 42 of the data table
# (pushing the former rows 42-147 down to 43-148), growing the sheet
# from A1:T147 to A1:T148.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 42; this shifts every
# row from 42 downward by one position (old row 42 -> new row 43, ...,
# old row 147 -> new row 148) and grows the sheet dimension to T148.
$ws.Rows(42).Insert()

# Populate the freshly inserted row 42 with the new record.
$ws.Range("A42").Value = 11
$ws.Range("B42").Value = "Vega Monumental Concepción"
$ws.Range("C42").Value = "Bíobío"
$ws.Range("D42").Value = 44622
$ws.Range("E42").Value = 8
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100108
$ws.Range("H42").Value = "Tropicales y subtropicales"
$ws.Range("I42").Value = 100108005
$ws.Range("J42").Value = "Piña"
$ws.Range("K42").Value = "Caramelo"
$ws.Range("L42").Value = "Segunda"
$ws.Range("M42").Value = 190
$ws.Range("N42").Value = 15000
$ws.Range("O42").Value = 16000
$ws.Range("P42").Value = 15474
$ws.Range("Q42").Value = "$/caja 14 unidades"
$ws.Range("R42").Value = "Ecuador"
$ws.Range("S42").Value = 1105
$ws.Range("T42").Value = 14
